# Break out stock.yaml completed:
#  - Fix bsecode column (D) for rows 567-572 from text to numeric values
#  - Append 11 new data rows (573-583) of "day" stock movers

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Convert bsecode (column D) on rows 567-572 from text to numbers ---
$bseFix = @{
    567 = 540719
    568 = 500043
    569 = 500770
    570 = 500112
    571 = 540777
    572 = 532483
}
foreach ($row in $bseFix.Keys) {
    $ws.Cells.Item($row, 4).Value = $bseFix[$row]
}

# --- 2. Append new rows 573-583 ---
$newRows = @(
    @(1,  "SHREECEM",   "Shree Cements Limited",                                    "500387", -0.7,  24900,    20392,    "day", "19/09/2024 11:35:35"),
    @(2,  "ULTRACEMCO", "Ultratech Cement Limited",                                  "532538", -0.03, 11627.5,  227247,   "day", "19/09/2024 11:35:35"),
    @(3,  "HDFCAMC",    "HDFC Asset Management Company Ltd",                         "541729", -0.96, 4389.8,   448620,   "day", "19/09/2024 11:35:35"),
    @(4,  "SRF",        "Srf Limited",                                               "503806", 0.57,  2402,     641706,   "day", "19/09/2024 11:35:35"),
    @(5,  "LUPIN",      "Lupin Limited",                                             "500257", -2.38, 2171.9,   1099230,  "day", "19/09/2024 11:35:35"),
    @(6,  "DALBHARAT",  "Dalmia Bharat Limited",                                     "533309", -0.67, 1822.25,  265631,   "day", "19/09/2024 11:35:35"),
    @(7,  "ICICIPRULI", "Icici Prudential Life Insurance Company Limited",           "540133", 0.59,  755.1,    1508414,  "day", "19/09/2024 11:35:35"),
    @(8,  "GNFC",       "Gujarat Narmada Valley Fertilizers And Chemicals Limited",  "500670", -1.32, 646.25,   963737,   "day", "19/09/2024 11:35:35"),
    @(9,  "PFC",        "Power Finance Corporation Limited",                         "532810", -2.31, 480.7,    15783626, "day", "19/09/2024 11:35:35"),
    @(10, "HINDPETRO",  "Hindustan Petroleum Corporation Limited",                   "500104", -2.31, 398.05,   4541762,  "day", "19/09/2024 11:35:35"),
    @(11, "IOC",        "Indian Oil Corporation Limited",                            "530965", -2.02, 165.04,   19200769, "day", "19/09/2024 11:35:35")
)

$startRow = 573
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]          # sr
    $ws.Cells.Item($r, 2).Value = $data[1]          # nsecode
    $ws.Cells.Item($r, 3).Value = $data[2]          # name

    # bsecode stays text (matches existing convention for not-yet-fixed rows);
    # leading apostrophe forces text entry instead of auto-numeric conversion
    $ws.Cells.Item($r, 4).Value = "'" + $data[3]    # bsecode

    $ws.Cells.Item($r, 5).Value = $data[4]          # per_chg
    $ws.Cells.Item($r, 6).Value = $data[5]          # close
    $ws.Cells.Item($r, 7).Value = $data[6]          # volume
    $ws.Cells.Item($r, 8).Value = $data[7]          # timeframe

    # Date Time stored as plain text, matching existing rows
    $ws.Cells.Item($r, 9).Value = $data[8]          # Date Time
}
